# Apply latest crypto price/volume updates scraped by GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.306.56"
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = "'3.310.93"
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'189.96"
$ws.Range('E5').Value = '  +2.98%  '
$ws.Range('D6').Value = "'562.28"
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.590"
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = "'3.303.36"
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').Value = "'0.186"
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').Value = "'0.589"
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = "'47.97"
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').Value = "'0.0000272"
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = "'8.74"
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').Value = "'3.843.00"
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').Value = "'616.08"
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = "'66.388.25"
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = "'18.13"
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = "'3.327.85"
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  -4.46%  '
$ws.Range('D22').Value = "'0.912"
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').Value = "'18.54"
$ws.Range('E23').Value = '  +8.83%  '
$ws.Range('D24').Value = "'5.12"
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('D25').Value = "'102.06"
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('D27').Value = "'6.02"
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('D29').Value = "'9.77"
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('D30').Value = "'8.64"
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('D31').Value = "'30.37"
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = "'6.76"
$ws.Range('E32').Value = '  +5.94%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D33').Value = "'4.11"
$ws.Range('E33').Value = '  +5.24%  '
$ws.Range('D34').Value = "'578.72"
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('D35').Value = "'11.14"
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').Value = "'3.746.56"
$ws.Range('E37').Value = '  -3.90%  '
$ws.Range('D38').Value = "'57.36"
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('D39').Value = "'0.999"
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').Value = "'3.57"
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('E42').Value = '  -3.90%  '
$ws.Range('D43').Value = "'34.16"
$ws.Range('E43').Value = '  +4.89%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = "'0.131"
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = "'2.75"
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('D46').Value = "'0.343"
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('D47').Value = "'0.0428"
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').Value = "'3.28"
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').Value = "'2.60"
$ws.Range('E50').Value = '  -3.41%  '
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  +0.15%  '
